# Cheetah Processing DEV test file refresh - 14th June 2022
#
# Refreshes the ShipmentTrackNum / PackageTrackNum columns with a newly
# generated batch of FedEx tracking numbers, and flips the PASS markers
# to FAIL for the scenario rows whose tracking numbers were refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ShipmentTrackNum values (column C). These are long digit-only
# strings, so the cell is pre-formatted as Text before the assignment
# (otherwise Excel auto-coerces a numeric-looking string into a Number)
# and the temporary Text number-format is cleared back to the workbook's
# Normal style immediately afterwards so no visible formatting changes.
function Set-TrackingNumber($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TrackingNumber "C2"  "310109786482"
Set-TrackingNumber "C3"  "310109786493"
Set-TrackingNumber "C4"  "310109786520"
Set-TrackingNumber "C5"  "310109786541"
Set-TrackingNumber "C6"  "310109786585"
Set-TrackingNumber "C7"  "310109786600"
Set-TrackingNumber "C8"  "310109786633"
Set-TrackingNumber "C9"  "310109786655"
Set-TrackingNumber "C10" "310109786688"
Set-TrackingNumber "C11" "310109786703"
Set-TrackingNumber "C12" "310109786747"
Set-TrackingNumber "C13" "310109786769"
Set-TrackingNumber "C14" "310109786791"
Set-TrackingNumber "C15" "310109786817"
Set-TrackingNumber "C16" "310109786840"
Set-TrackingNumber "C17" "310109786861"
Set-TrackingNumber "C18" "310109786909"
Set-TrackingNumber "C19" "310109786920"
Set-TrackingNumber "C20" "310109786953"
Set-TrackingNumber "C21" "310109786975"

# New PackageTrackNum values (column D) - mirrors column C on these rows
Set-TrackingNumber "D5"  "310109786541"
Set-TrackingNumber "D6"  "310109786585"
Set-TrackingNumber "D7"  "310109786600"
Set-TrackingNumber "D13" "310109786769"
Set-TrackingNumber "D14" "310109786791"
Set-TrackingNumber "D15" "310109786817"
Set-TrackingNumber "D16" "310109786840"
Set-TrackingNumber "D17" "310109786861"

# Flip PASS -> FAIL for the status columns tied to the refreshed rows
$failCells = @(
    "L2","M2","N2","O2","P2",
    "Q3",
    "M4","N4","O4","P4","R4",
    "M5","N5","O5","P5","S5",
    "L6","M6","N6","O6","P6","S6",
    "L7","M7","N7","O7","P7","T7",
    "L13","M13","N13","O13","P13","U13",
    "L14","M14","N14","O14","P14","V14",
    "L15","M15","N15","O15","P15","W15",
    "L16","M16","N16","O16","P16","X16",
    "L17","M17","N17","O17","P17","Y17",
    "L19","M19","N19","O19","P19","AB19",
    "L20","M20","N20","O20","P20","AC20",
    "L21","M21","N21","O21","P21","AD21"
)

foreach ($addr in $failCells) {
    $ws.Range($addr).Value = "FAIL"
}
